# Bugfixed evaluation and simulated rt_data for components
# Insert a new row of data (2007 -> 2008) at the top of the data block,
# shifting everything down by one row, and populate the previously
# missing y_0_forecast (column C) value for that first data row, plus
# recompute all subsequent rows' values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 2 so the existing data (rows 2:18) shifts
# down to rows 3:19.
$ws.Rows.Item(2).Insert(-4121) # xlShiftDown

# The freshly inserted row inherits the formatting of the row above it
# (the bold/bordered header row). Reset it back to the plain/default
# style used by the rest of the data rows, then re-apply the date
# number format (style used by column A in the other data rows) only
# to the new A2 cell.
$ws.Rows.Item(2).ClearFormats()
$ws.Range("A3").Copy()
$ws.Range("A2").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

$data = @(
    @(39400, 2007, 4.930115226412335, 2008, 0.3630458632513767),
    @(39765, 2008, 1.457587285166628, 2009, 0.507956838644974),
    @(40130, 2009, -0.9140166223623569, 2010, 5.6395352704941),
    @(40494, 2010, 2.585942866987878, 2011, 1.724360951547554),
    @(40862, 2011, 4.253963781362402, 2012, -0.6955733540840336),
    @(41228, 2012, 1.752870900283909, 2013, 4.300339264728548),
    @(41592, 2013, -1.479696720105139, 2014, 8.296896928314457),
    @(41957, 2014, 3.900127535411246, 2015, -2.092856741436233),
    @(42321, 2015, 0.03947433952959933, 2016, -1.259568900987018),
    @(42689, 2016, 2.192778679161944, 2017, 0.5033587260849126),
    @(43053, 2017, 3.40836448860673, 2018, -0.3010260522302355),
    @(43418, 2018, 2.799070570134488, 2019, 3.825329033908775),
    @(43783, 2019, 4.195393191694419, 2020, 3.942709467505678),
    @(44159, 2020, 1.666553973046048, 2021, -0.2638638106667313),
    @(44525, 2021, 1.879266440112803, 2022, -0.04532879466145889),
    @(44890, 2022, -2.620683231370946, 2023, -3.179374983142691),
    @(45254, 2023, -3.036556262700274, 2024, -2.51939929628594),
    @(45618, 2024, -2.953443685011514, 2025, -2.112604539331953)
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = $r[4]
    $row++
}
